$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.887.37'
$ws.Range('E2').Value = '  +4.26%  '
$ws.Range('D3').Value = '2.779.21'
$ws.Range('E3').Value = '  +4.56%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '343.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '115.38'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.549'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.577'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.70'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0855'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.02'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('D15').Value = '3.210.25'
$ws.Range('E15').Value = '  +4.43%  '
$ws.Range('D16').Value = '2.777.54'
$ws.Range('E16').Value = '  +4.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.882'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.74%  '
$ws.Range('D18').Value = '51.810.63'
$ws.Range('E18').Value = '  +4.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.90%  '
$ws.Range('E20').Value = '  +4.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.98%  '
$ws.Range('D22').Value = '0.0₃0979'
$ws.Range('E22').Value = '  +2.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '271.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('E25').Value = '  +7.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.22'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.57'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.72'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0820'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.10'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.57%  '
$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.94'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.23'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0385'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.66'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +24.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '127.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '23.35'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.63%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.33'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').Value = '2.070.31'
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.52%  '
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.55'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.900'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.91%  '
$ws.Range('E51').Value = '  -1.63%  '
